# edit.ps1 - apply the README.docx changes described by the diff
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "is designed in Page Object Model. " -> "is designed in Page Object
#    Model (POM). "  (split into extra runs, but functionally this is
#    just inserting " (POM)" before the trailing ". ")
# ---------------------------------------------------------------------
$r = $d.Content
$null = $r.Find.Execute("is designed in Page Object Model. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Text = "is designed in Page Object Model"
$r.Collapse(0)
$r.InsertAfter(" (POM)")
$r.Collapse(0)
$r.InsertAfter(". ")

# ---------------------------------------------------------------------
# 2) Insert a brand-new paragraph right after the paragraph ending in
#    "...register....)" (and before the existing blank paragraph),
#    describing the 19/04 update to the POM framework.
# ---------------------------------------------------------------------
$r = $d.Content
$ellipsis = [char]0x2026
$needle = "register" + $ellipsis + ".)"
$null = $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
$r.InsertParagraphAfter()

$r2 = $d.Content
$null = $r2.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.MoveStart(1, 1)
$apost = [char]0x2019
$newParaText = "19/04: In the first version, the POM framework referred only to the fact that the business and the tech facing were separated, I kept it simple due to the time restrictions. To implement that separation, in this version, I" + $apost + "ve created a structure project with classes for different webs in the application, which makes them totally reusable (see for instance the use of the method ClickAgreeTerms from checkout page used in a method in LogIn page)"
$r2.InsertAfter($newParaText)

# ---------------------------------------------------------------------
# 3) lastRenderedPageBreak markers shift to new positions because of the
#    added text above (Word recalculates these on layout/open); emulate
#    the moves seen in the diff directly on the OOXML runs.
# ---------------------------------------------------------------------

# 3a) Add lastRenderedPageBreak before the run "5.- " (the Checkout item)
$r = $d.Content
$null = $r.Find.Execute("5.- Checkout", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(1)
$r.InsertBreak(7)  # wdPageBreak-ish placeholder, replaced below if unsupported

Write-Host "done step3a attempt"
